$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (Colaborador_id, Colaborador_nome, Departamento,
# Motivo_da_ausência, Horas_de_ausência, Data_da_ausência, Salário)
$data = @(
    @(72085, "Sabrina Peixoto", "Operações", "Problemas pessoais", 1, 45081, 8909.540000000001),
    @(45684, "Dr. Joaquim Aragão", "Atendimento ao Cliente", "Problemas pessoais", 6, 45099, 8774.48),
    @(64045, "Cecília Almeida", "Marketing", "Viagem de negócios", 1, 45078, 3473.42),
    @(26136, "Emanuelly Caldeira", "Engenharia", "Doença", 4, 45081, 5187.16),
    @(85292, "Milena Araújo", "Atendimento ao Cliente", "Consulta médica", 7, 45102, 7591.83),
    @(85935, "Paulo Araújo", "Marketing", "Doença", 6, 45086, 9084.84),
    @(83712, "Isabella Rezende", "Engenharia", "Doença", 7, 45094, 4485.87),
    @(27544, "Sophia Martins", "Engenharia", "Problemas pessoais", 3, 45106, 2793.64),
    @(78630, "Eduarda Campos", "P&D", "Doença", 3, 45092, 6834.95),
    @(96749, "Amanda Melo", "Jurídico", "Viagem de negócios", 2, 45102, 8941.16)
)

$rowIndex = 2
foreach ($rowData in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowData[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rowData[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rowData[5]
    $ws.Cells.Item($rowIndex, 7).Value = $rowData[6]
    $rowIndex++
}
